$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.631305391513564
$ws.Range("C2").Value = 0.210090626371759
$ws.Range("D2").Value = 0.07955070600113245
$ws.Range("F2").Value = 1.671718243172151
$ws.Range("G2").Value = 1.597995458338389
$ws.Range("H2").Value = 1.417263889856883
$ws.Range("K2").Value = 0.2753409718590092
$ws.Range("L2").Value = 0.2856578238929046
$ws.Range("M2").Value = 0.2063240920548068
$ws.Range("B3").Value = 0.6009864536892167
$ws.Range("C3").Value = 0.209270069445143
$ws.Range("D3").Value = 0.07924658692274278
$ws.Range("F3").Value = 1.656374741911378
$ws.Range("G3").Value = 1.58291719738564
$ws.Range("H3").Value = 1.414682555761416
$ws.Range("K3").Value = 0.2479803965992744
$ws.Range("L3").Value = 0.2821316541934635
$ws.Range("M3").Value = 0.2000368177901315
$ws.Range("B4").Value = 0.5827447636306431
$ws.Range("C4").Value = 0.2087463781081382
$ws.Range("D4").Value = 0.07904912578464618
$ws.Range("F4").Value = 1.647677632131206
$ws.Range("G4").Value = 1.574351356434406
$ws.Range("H4").Value = 1.413585509933895
$ws.Range("K4").Value = 0.2312963188681465
$ws.Range("L4").Value = 0.28010393415029
$ws.Range("M4").Value = 0.1962942113172019
$ws.Range("B5").Value = 0.5754054988939572
$ws.Range("C5").Value = 0.2085279639586695
$ws.Range("D5").Value = 0.07896596391388755
$ws.Range("F5").Value = 1.644315427856313
$ws.Range("G5").Value = 1.571034606642399
$ws.Range("H5").Value = 1.413261169919224
$ws.Range("K5").Value = 0.2245265738137618
$ws.Range("L5").Value = 0.2793121947027473
$ws.Range("M5").Value = 0.1947987466454961
$ws.Range("B6").Value = 0.5741925272065203
$ws.Range("C6").Value = 0.2084913938271313
$ws.Range("D6").Value = 0.078951992291179
$ws.Range("F6").Value = 1.643768122671872
$ws.Range("G6").Value = 1.570494360832285
$ws.Range("H6").Value = 1.413214726036102
$ws.Range("K6").Value = 0.2234042273504855
$ws.Range("L6").Value = 0.2791828160484968
$ws.Range("M6").Value = 0.1945522199783305
$ws.Range("B7").Value = 0.5826454012814963
$ws.Range("C7").Value = 0.2087434527798742
$ws.Range("D7").Value = 0.07904801514133553
$ws.Range("F7").Value = 1.647631551721048
$ws.Range("G7").Value = 1.574305921712948
$ws.Range("H7").Value = 1.413580638871906
$ws.Range("K7").Value = 0.2312049015118873
$ws.Range("L7").Value = 0.2800931164476026
$ws.Range("M7").Value = 0.1962739227197403
$ws.Range("B8").Value = 0.6207738885909748
$ws.Range("C8").Value = 0.2098118152638335
$ws.Range("D8").Value = 0.07944807493150918
$ws.Range("F8").Value = 1.666277508556405
$ws.Range("G8").Value = 1.592652678298904
$ws.Range("H8").Value = 1.416272594573527
$ws.Range("K8").Value = 0.2658831554872165
$ws.Range("L8").Value = 0.2844135083424817
$ws.Range("M8").Value = 0.2041318195078574
$ws.Range("B9").Value = 0.6985066554974253
$ws.Range("C9").Value = 0.2117497513186954
$ws.Range("D9").Value = 0.08014731183374302
$ws.Range("F9").Value = 1.708593051412393
$ws.Range("G9").Value = 1.634134928564094
$ws.Range("H9").Value = 1.425423141294203
$ws.Range("K9").Value = 0.3348011745949293
$ws.Range("L9").Value = 0.2939750740664095
$ws.Range("M9").Value = 0.2204745537251185
$ws.Range("B10").Value = 0.7574213457708368
$ws.Range("C10").Value = 0.2130785428536832
$ws.Range("D10").Value = 0.0806088977738213
$ws.Range("F10").Value = 1.743203611064132
$ws.Range("G10").Value = 1.667988324927848
$ws.Range("H10").Value = 1.434509046435807
$ws.Range("K10").Value = 0.3859957132436307
$ws.Range("L10").Value = 0.3016643543711837
$ws.Range("M10").Value = 0.2330504597391396
$ws.Range("B11").Value = 0.7846150875613773
$ws.Range("C11").Value = 0.2136625847194118
$ws.Range("D11").Value = 0.08080753351047676
$ws.Range("F11").Value = 1.759717123631447
$ws.Range("G11").Value = 1.684127041033321
$ws.Range("H11").Value = 1.439156315705333
$ws.Range("K11").Value = 0.4094082126240437
$ws.Range("L11").Value = 0.3053068287745759
$ws.Range("M11").Value = 0.2388951803740937
$ws.Range("B12").Value = 0.7949690433895853
$ws.Range("C12").Value = 0.2138808200552837
$ws.Range("D12").Value = 0.08088111794493358
$ws.Range("F12").Value = 1.766081127915058
$ws.Range("G12").Value = 1.690344850133471
$ws.Range("H12").Value = 1.440990053792291
$ws.Range("K12").Value = 0.4182916931493992
$ws.Range("L12").Value = 0.3067069186187439
$ws.Range("M12").Value = 0.2411262104547944
$ws.Range("B13").Value = 0.7927366349472322
$ws.Range("C13").Value = 0.2138339492720007
$ws.Range("D13").Value = 0.08086534297637371
$ws.Range("F13").Value = 1.764705600688202
$ws.Range("G13").Value = 1.689000996917827
$ws.Range("H13").Value = 1.440591837969123
$ws.Range("K13").Value = 0.4163776926729668
$ws.Range("L13").Value = 0.3064044611680856
$ws.Range("M13").Value = 0.2406449290993038
$ws.Range("B14").Value = 0.7854657873561166
$ws.Range("C14").Value = 0.2136805977120346
$ws.Range("D14").Value = 0.08081362012708126
$ws.Range("F14").Value = 1.760238474935619
$ws.Range("G14").Value = 1.684636449747899
$ws.Range("H14").Value = 1.439305697017716
$ws.Range("K14").Value = 0.4101387084361932
$ws.Range("L14").Value = 0.3054215990236315
$ws.Range("M14").Value = 0.2390783728973886
$ws.Range("B15").Value = 0.781019504428059
$ws.Range("C15").Value = 0.2135862843404297
$ws.Range("D15").Value = 0.08078172540856343
$ws.Range("F15").Value = 1.757516650646053
$ws.Range("G15").Value = 1.681976904892082
$ws.Range("H15").Value = 1.438527524840538
$ws.Range("K15").Value = 0.4063194483341306
$ws.Range("L15").Value = 0.3048222708519575
$ws.Range("M15").Value = 0.2381211234851648
$ws.Range("B16").Value = 0.7556520598341478
$ws.Range("C16").Value = 0.2130399640057981
$ws.Range("D16").Value = 0.08059568791563265
$ws.Range("F16").Value = 1.742139898161255
$ws.Range("G16").Value = 1.666948503309555
$ws.Range("H16").Value = 1.434215680268125
$ws.Range("K16").Value = 0.3844681319629331
$ws.Range("L16").Value = 0.3014292171593524
$ws.Range("M16").Value = 0.2326709817582966
$ws.Range("B17").Value = 0.7401904650193387
$ws.Range("C17").Value = 0.2126995877424953
$ws.Range("D17").Value = 0.08047865233906393
$ws.Range("F17").Value = 1.732903798581816
$ws.Range("G17").Value = 1.657918399619916
$ws.Range("H17").Value = 1.431702161773813
$ws.Range("K17").Value = 0.3710946812644238
$ws.Range("L17").Value = 0.2993847003494778
$ws.Range("M17").Value = 0.2293591870579945
$ws.Range("B18").Value = 0.7313343854002596
$ws.Range("C18").Value = 0.2125018885870773
$ws.Range("D18").Value = 0.08041026921266337
$ws.Range("F18").Value = 1.727663813013464
$ws.Range("G18").Value = 1.652794034074645
$ws.Range("H18").Value = 1.430304836817612
$ws.Range("K18").Value = 0.3634142910571256
$ws.Range("L18").Value = 0.2982223587110724
$ws.Range("M18").Value = 0.227465993402916
$ws.Range("B19").Value = 0.7283422353069966
$ws.Range("C19").Value = 0.2124346203594811
$ws.Range("D19").Value = 0.08038693265519825
$ws.Range("F19").Value = 1.725902070848235
$ws.Range("G19").Value = 1.651070944095778
$ws.Range("H19").Value = 1.429840036711227
$ws.Range("K19").Value = 0.3608158504804351
$ws.Range("L19").Value = 0.2978311482201548
$ws.Range("M19").Value = 0.2268269953056432
$ws.Range("B20").Value = 0.7418325485584205
$ws.Range("C20").Value = 0.2127360204084319
$ws.Range("D20").Value = 0.08049122147006571
$ws.Range("F20").Value = 1.733879505809043
$ws.Range("G20").Value = 1.658872473873146
$ws.Range("H20").Value = 1.431964722507956
$ws.Range("K20").Value = 0.3725171016255615
$ws.Range("L20").Value = 0.2996009342306252
$ws.Range("M20").Value = 0.2297105265553654
$ws.Range("B21").Value = 0.7875998862765812
$ws.Range("C21").Value = 0.2137257201888012
$ws.Range("D21").Value = 0.08082885677268337
$ws.Range("F21").Value = 1.761547572896731
$ws.Range("G21").Value = 1.685915532944904
$ws.Range("H21").Value = 1.439681461681687
$ws.Range("K21").Value = 0.4119707706989288
$ws.Range("L21").Value = 0.3057097261485637
$ws.Range("M21").Value = 0.2395380268067981
$ws.Range("B22").Value = 0.8178393887812945
$ws.Range("C22").Value = 0.2143554802694965
$ws.Range("D22").Value = 0.0810399916289235
$ws.Range("F22").Value = 1.780275525149534
$ws.Range("G22").Value = 1.704210174945786
$ws.Range("H22").Value = 1.445155647805194
$ws.Range("K22").Value = 0.4378589879165986
$ws.Range("L22").Value = 0.3098231900933399
$ws.Range("M22").Value = 0.2460643685142401
$ws.Range("B23").Value = 0.8016700789367235
$ws.Range("C23").Value = 0.214020923852388
$ws.Range("D23").Value = 0.08092817808124408
$ws.Range("F23").Value = 1.770220984371846
$ws.Range("G23").Value = 1.694389139849363
$ws.Range("H23").Value = 1.442194548974157
$ws.Range("K23").Value = 0.4240325816005566
$ws.Range("L23").Value = 0.3076166928699422
$ws.Range("M23").Value = 0.2425716837716365
$ws.Range("B24").Value = 0.7410900597897694
$ws.Range("C24").Value = 0.2127195554698247
$ws.Range("D24").Value = 0.08048554238388306
$ws.Range("F24").Value = 1.733438170616594
$ws.Range("G24").Value = 1.658440927709279
$ws.Range("H24").Value = 1.43184587011271
$ws.Range("K24").Value = 0.3718739998783747
$ws.Range("L24").Value = 0.2995031341546763
$ws.Range("M24").Value = 0.2295516523207723
$ws.Range("B25").Value = 0.6771608851508404
$ws.Range("C25").Value = 0.2112422292287306
$ws.Range("D25").Value = 0.07996729604591302
$ws.Range("F25").Value = 1.696528307918086
$ws.Range("G25").Value = 1.622321453241852
$ws.Range("H25").Value = 1.422532909225225
$ws.Range("K25").Value = 0.3160588428297615
$ws.Range("L25").Value = 0.2912717675868919
$ws.Range("M25").Value = 0.2784865382909321
